$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 31 de Marzo de 2020 a las 20:50'

$ws.Range("B4").Value = 180340
$ws.Range("C4").Value = 16552
$ws.Range("E4").Value = 170525
$ws.Range("F4").Value = 3981
$ws.Range("G4").Value = 433
$ws.Range("H4").Value = 3574

$ws.Range("A31").Value = 'Polonia'
$ws.Range("B31").Value = 2311
$ws.Range("C31").Value = 256
$ws.Range("D31").Value = 7
$ws.Range("E31").Value = 2271
$ws.Range("F31").Value = 50
$ws.Range("G31").Value = 2
$ws.Range("H31").Value = 33

$ws.Range("A32").Value = 'Rumania'
$ws.Range("B32").Value = 2245
$ws.Range("C32").Value = 136
$ws.Range("D32").Value = 220
$ws.Range("E32").Value = 1945
$ws.Range("F32").Value = 62
$ws.Range("G32").Value = 15
$ws.Range("H32").Value = 80

$ws.Range("A33").Value = 'Ecuador'
$ws.Range("B33").Value = 2240
$ws.Range("C33").Value = 274
$ws.Range("D33").Value = 54
$ws.Range("E33").Value = 2111
$ws.Range("F33").Value = 100
$ws.Range("G33").Value = 13
$ws.Range("H33").Value = 75

$ws.Range("B37").Value = 1938
$ws.Range("C37").Value = 221
$ws.Range("E37").Value = 1836

$ws.Range("A42").Value = 'India'
$ws.Range("B42").Value = 1397
$ws.Range("C42").Value = 146
$ws.Range("D42").Value = 123
$ws.Range("E42").Value = 1239
$ws.Range("F42").Value = 0
$ws.Range("G42").Value = 3
$ws.Range("H42").Value = 35

$ws.Range("A43").Value = 'Sudafrica'
$ws.Range("B43").Value = 1353
$ws.Range("C43").Value = 27
$ws.Range("D43").Value = 31
$ws.Range("E43").Value = 1319
$ws.Range("F43").Value = 7
$ws.Range("G43").Value = 0
$ws.Range("H43").Value = 3

$ws.Range("A44").Value = 'Grecia'
$ws.Range("B44").Value = 1314
$ws.Range("C44").Value = 102
$ws.Range("D44").Value = 52
$ws.Range("E44").Value = 1213
$ws.Range("F44").Value = 72
$ws.Range("G44").Value = 3
$ws.Range("H44").Value = 49

$ws.Range("D49").Value = 394
$ws.Range("E49").Value = 647

$ws.Range("A56").Value = 'Catar'
$ws.Range("B56").Value = 781
$ws.Range("C56").Value = 88
$ws.Range("D56").Value = 62
$ws.Range("E56").Value = 717
$ws.Range("F56").Value = 6
$ws.Range("H56").Value = 2

$ws.Range("A57").Value = 'Estonia'
$ws.Range("B57").Value = 745
$ws.Range("C57").Value = 30
$ws.Range("D57").Value = 26
$ws.Range("E57").Value = 715
$ws.Range("F57").Value = 13
$ws.Range("G57").Value = 1
$ws.Range("H57").Value = 4

$ws.Range("A58").Value = 'Argelia'
$ws.Range("B58").Value = 716
$ws.Range("C58").Value = 132
$ws.Range("D58").Value = 46
$ws.Range("E58").Value = 626
$ws.Range("F58").Value = 0
$ws.Range("G58").Value = 9
$ws.Range("H58").Value = 44

$ws.Range("A59").Value = 'Hong Kong'
$ws.Range("B59").Value = 714
$ws.Range("C59").Value = 31
$ws.Range("D59").Value = 128
$ws.Range("E59").Value = 582
$ws.Range("F59").Value = 5
$ws.Range("H59").Value = 4

$ws.Range("A60").Value = 'Crucero'
$ws.Range("B60").Value = 712
$ws.Range("C60").Value = 0
$ws.Range("D60").Value = 603
$ws.Range("E60").Value = 99
$ws.Range("F60").Value = 15
$ws.Range("G60").Value = 0
$ws.Range("H60").Value = 10

$ws.Range("A61").Value = 'Egipto'
$ws.Range("B61").Value = 710
$ws.Range("C61").Value = 54
$ws.Range("D61").Value = 157
$ws.Range("E61").Value = 507
$ws.Range("F61").Value = 0
$ws.Range("G61").Value = 5
$ws.Range("H61").Value = 46

$ws.Range("A62").Value = 'Irak'
$ws.Range("B62").Value = 694
$ws.Range("C62").Value = 64
$ws.Range("D62").Value = 170
$ws.Range("E62").Value = 474
$ws.Range("F62").Value = 0
$ws.Range("G62").Value = 4
$ws.Range("H62").Value = 50

$ws.Range("A63").Value = 'Emiratos Arabes Unidos'
$ws.Range("B63").Value = 664
$ws.Range("C63").Value = 53
$ws.Range("D63").Value = 61
$ws.Range("E63").Value = 597
$ws.Range("F63").Value = 2
$ws.Range("G63").Value = 1
$ws.Range("H63").Value = 6

$ws.Range("E104").Value = 138
$ws.Range("G104").Value = 2
$ws.Range("H104").Value = 5

$ws.Range("B110").Value = 128
$ws.Range("C110").Value = 35
$ws.Range("E110").Value = 98
$ws.Range("G110").Value = 2
$ws.Range("H110").Value = 3

$ws.Range("A112").Value = 'Guadalupe'
$ws.Range("B112").Value = 114
$ws.Range("C112").Value = 8
$ws.Range("D112").Value = 22
$ws.Range("E112").Value = 88
$ws.Range("F112").Value = 14
$ws.Range("H112").Value = 4

$ws.Range("A113").Value = 'Georgia'
$ws.Range("B113").Value = 110
$ws.Range("C113").Value = 7
$ws.Range("D113").Value = 21
$ws.Range("E113").Value = 89
$ws.Range("F113").Value = 6
$ws.Range("G113").Value = 0
$ws.Range("H113").Value = 0

$ws.Range("A114").Value = 'Montenegro'
$ws.Range("C114").Value = 18
$ws.Range("D114").Value = 0
$ws.Range("E114").Value = 107
$ws.Range("G114").Value = 1
$ws.Range("H114").Value = 2

$ws.Range("A115").Value = 'Camboya'
$ws.Range("B115").Value = 109
$ws.Range("C115").Value = 2
$ws.Range("D115").Value = 23
$ws.Range("E115").Value = 86
$ws.Range("F115").Value = 1

$ws.Range("A116").Value = 'Kirguistan'
$ws.Range("C116").Value = 13
$ws.Range("D116").Value = 3
$ws.Range("E116").Value = 104
$ws.Range("F116").Value = 0
$ws.Range("G116").Value = 0
$ws.Range("H116").Value = 0

$ws.Range("A117").Value = 'Bolivia'
$ws.Range("B117").Value = 107
$ws.Range("C117").Value = 10
$ws.Range("D117").Value = 0
$ws.Range("E117").Value = 101
$ws.Range("F117").Value = 3
$ws.Range("G117").Value = 2
$ws.Range("H117").Value = 6

$ws.Range("B121").Value = 75
$ws.Range("C121").Value = 5
$ws.Range("E121").Value = 75

$ws.Range("A127").Value = 'Madagascar'
$ws.Range("B127").Value = 57
$ws.Range("C127").Value = 14
$ws.Range("D127").Value = 0
$ws.Range("E127").Value = 57
$ws.Range("F127").Value = 6

$ws.Range("A128").Value = 'Aruba'
$ws.Range("B128").Value = 55
$ws.Range("C128").Value = 5
$ws.Range("D128").Value = 1
$ws.Range("E128").Value = 54
$ws.Range("H128").Value = 0

$ws.Range("A129").Value = 'Monaco'
$ws.Range("B129").Value = 52
$ws.Range("C129").Value = 3
$ws.Range("D129").Value = 2
$ws.Range("E129").Value = 49
$ws.Range("F129").Value = 0
$ws.Range("H129").Value = 1

$ws.Range("A130").Value = 'Banglades'
$ws.Range("B130").Value = 51
$ws.Range("C130").Value = 2
$ws.Range("D130").Value = 25
$ws.Range("E130").Value = 21
$ws.Range("F130").Value = 1
$ws.Range("H130").Value = 5
